$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats/styles from column E into the new column D
# (restricted to the row ranges that actually contain data, to avoid
# touching blank separator rows 5,6,36,37,78,79)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 706100
$ws.Range("D9").Value = 449700
$ws.Range("D10").Value = 256400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 85200
$ws.Range("D17").Value = 574000
$ws.Range("D18").Value = 132200
$ws.Range("D20").Value = -3600
$ws.Range("D21").Value = 213700
$ws.Range("D22").Value = 37100
$ws.Range("D23").Value = 91500
$ws.Range("D24").Value = 24200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 67300
$ws.Range("D27").Value = 67300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -2700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 3600
$ws.Range("D33").Value = 64600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 64600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 12600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 72000
$ws.Range("D44").Value = 44100
$ws.Range("D45").Value = 167100
$ws.Range("D46").Value = 295900
$ws.Range("D47").Value = 63600
$ws.Range("D48").Value = 2421400
$ws.Range("D49").Value = 9000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 452900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3242700
$ws.Range("D57").Value = 115900
$ws.Range("D58").Value = 247600
$ws.Range("D59").Value = 145600
$ws.Range("D60").Value = 509100
$ws.Range("D61").Value = 706200
$ws.Range("D62").Value = 1264700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2480000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 312200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 762600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 64600
$ws.Range("D83").Value = 85200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 168800
$ws.Range("D91").Value = -214600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -217500
$ws.Range("D96").Value = -51300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 57800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 9200

# Apply the individual figure corrections/updates in the shifted columns
$ws.Range("E8").Value = 755000
$ws.Range("F8").Value = 668200
$ws.Range("E9").Value = 469400
$ws.Range("F9").Value = 631000
$ws.Range("E10").Value = 285700
$ws.Range("F10").Value = 37200
$ws.Range("E15").Value = 81100
$ws.Range("F15").Value = 77600
$ws.Range("E17").Value = 604100
$ws.Range("F17").Value = 517500
$ws.Range("E18").Value = 150900
$ws.Range("F18").Value = 150700
$ws.Range("E20").Value = -300
$ws.Range("F20").Value = -7200
$ws.Range("E21").Value = 231700
$ws.Range("F21").Value = 221200
$ws.Range("E22").Value = 37500
$ws.Range("F22").Value = 38100
$ws.Range("E23").Value = 113100
$ws.Range("F23").Value = 105400
$ws.Range("E24").Value = 44400
$ws.Range("F24").Value = 43000
$ws.Range("E26").Value = 68700
$ws.Range("F26").Value = 62400
$ws.Range("E27").Value = 68700
$ws.Range("F27").Value = 62400
$ws.Range("E29").Value = -124300
$ws.Range("F29").Value = -3500
$ws.Range("E32").Value = 300
$ws.Range("F32").Value = 7200
$ws.Range("E43").Value = 65300
$ws.Range("E44").Value = 47600
$ws.Range("E45").Value = 153600
$ws.Range("E48").Value = 2244200
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "NA"
$ws.Range("H49").Value = "NA"
$ws.Range("I49").Value = "NA"
$ws.Range("J49").Value = "NA"
$ws.Range("E52").Value = 459300
$ws.Range("E57").Value = 111000
$ws.Range("E59").Value = 119900
$ws.Range("E62").Value = 1219900
$ws.Range("E83").Value = 81100
$ws.Range("F83").Value = 77600
$ws.Range("E91").Value = -213300
$ws.Range("F91").Value = -138400
